# Edit for "SparkApplication/EnrichedLog Format.xlsx":
#   - Rename the last field from "Timestamp_ingestion" to "TimestampIngestion"
#   - Clarify its description to mention the unit (milliseconds)
#   - Move the active selection to C39 (the Description cell of that row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EnrichedLog Format")

$ws.Range("A39").Value = "TimestampIngestion"
$ws.Range("C39").Value = "timestamp(in millisec) when the logs are parsed and inserted in the table. Needed to manage properly the analysis step"

$ws.Activate()
$ws.Range("C39").Select()
